$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    3  = @{ D = 117638.3569498936;  E = -0.06453857968723657; F = 0.199131106982574;  G = -1.193211556734104;  H = 13.53117379388421 }
    5  = @{ D = 119152.4858036125;  E = -0.05549820509275658; F = 0.2376122136663193; G = -1.175633438585491;  H = 12.45594312940294 }
    6  = @{ D = 119899.3395965008;  E = -0.05393258646426717; F = 0.2366695189903922; G = -0.8389677413715874; H = 8.972387853633752 }
    7  = @{ D = 120592.2997283303;  E = -0.05798408101535789; F = 0.2611744526057739; G = -1.199080725176566;  H = 10.69309257821405 }
    8  = @{ D = 122071.8969801087;  E = -0.06355390454484852; F = 0.20867344261672;   G = -0.7201054825622363; H = 6.613230798559294 }
    9  = @{ D = 123539.4464871183;  E = -0.09842058364807636; F = 0.3414427309120105; G = -1.594008665007883;  H = 9.835577155593709 }
    10 = @{ D = 125108.7987202931;  E = -0.1336230559611806;  F = 0.4506868660390958; G = -1.949595690875405;  H = 9.842349855930896 }
    11 = @{ D = 127188.1808784983;  E = -0.2293438851710582;  F = 0.9551295382658115; G = -3.08624083731187;   H = 16.36436107875861 }
    17 = @{ D = 117673.6121162518;  E = -0.07670429875435213; F = 0.1780960075230125; G = -0.7776599544290737; H = 8.753653927209092 }
    20 = @{ D = 117888.3684726979;  E = -0.05683847252392848; F = 0.1754336564964335; G = -0.124369869136962;  H = 5.653229328298232 }
}

foreach ($row in $updates.Keys) {
    $cols = $updates[$row]
    $ws.Range("D$row").Value = $cols.D
    $ws.Range("E$row").Value = $cols.E
    $ws.Range("F$row").Value = $cols.F
    $ws.Range("G$row").Value = $cols.G
    $ws.Range("H$row").Value = $cols.H
}
